$wb = $excel.ActiveWorkbook

# --- Sheet 1: Категории ---
$ws1 = $wb.Worksheets.Item("Категории")
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "Люкс"
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Премимум"
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Президентский"
$ws1.Columns.Item(2).ColumnWidth = 14.307291666666666

# --- Sheet 2: Номера ---
$ws2 = $wb.Worksheets.Item("Номера")
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 155
$ws2.Range("D2").Value = 2
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 156
$ws2.Range("D3").Value = 4
$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = 175
$ws2.Range("D4").Value = 2

# --- Sheet 3: Граждане ---
$ws3 = $wb.Worksheets.Item("Граждане")
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "Пупкин"
$ws3.Range("C2").Value = 4919234567
$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = "Бэкинн"
$ws3.Range("C3").Value = 4919234567
$ws3.Range("A4").Value = 3
$ws3.Range("B4").Value = "Мартирасян"
$ws3.Range("C4").Value = 4919234567
$ws3.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws3.Columns.Item(3).ColumnWidth = 10.166666666666666

# --- Sheet 4: Размещение ---
$ws4 = $wb.Worksheets.Item("Размещение")
$ws4.Range("B2").Value = 1
$ws4.Range("C2").Value = 1
$ws4.Range("D2").NumberFormat = "mm-dd-yy"
$ws4.Range("D2").Value = 44958
$ws4.Range("E2").Value = "7 дней"
$ws4.Range("B3").Value = 2
$ws4.Range("C3").Value = 2
$ws4.Range("D3").Value = 45140
$ws4.Range("E3").Value = "2дня"
$ws4.Range("B4").Value = 3
$ws4.Range("C4").Value = 3
$ws4.Range("D4").Value = 45204
$ws4.Range("E4").Value = "3 дня"

# Reuse D2's date style for D3:D4 instead of allocating duplicate style entries
$ws4.Range("D2").Copy()
$ws4.Range("D3:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selections (set last so final selection per-sheet sticks) ---
$ws1.Columns.Item(2).Select()
$ws2.Range("D4").Select()
$ws3.Range("C4").Select()
$ws4.Range("F4").Select()
